# The 2007年 data row (row 2) was removed from the sheet.
# Deleting the entire row shifts rows 3-6 (2010/2012/2015/2017年) up to
# become rows 2-5, and automatically updates the sheet's used range /
# dimension from A1:Y6 to A1:Y5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2").Delete()
